$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 6435.2085
$ws.Range("I115").Value = 551.36365
$ws.Range("J115").Value = 11413.846
$ws.Range("K115").Value = 1654.09095
$ws.Range("L115").Value = 34241.538
$ws.Range("M115").Value = -87.09095000000002
$ws.Range("N115").Value = -37375.538

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3168915.5
$ws.Range("I132").Value = 866646.4
$ws.Range("J132").Value = 7937901.5
$ws.Range("K132").Value = 2599939.2
$ws.Range("L132").Value = 23813704.5
$ws.Range("M132").Value = -2597409.2
$ws.Range("N132").Value = -23818764.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 27796660
$ws.Range("I137").Value = 6579720
$ws.Range("J137").Value = 78186900
$ws.Range("K137").Value = 19739160
$ws.Range("L137").Value = 234560700
$ws.Range("M137").Value = -19736610
$ws.Range("N137").Value = -234565800

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3257.9812
$ws.Range("I138").Value = 2920.6086
$ws.Range("J138").Value = 3516.6333
$ws.Range("K138").Value = 8761.825800000001
$ws.Range("L138").Value = 10549.8999
$ws.Range("M138").Value = -3621.825800000001
$ws.Range("N138").Value = -20829.8999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16975.137
$ws.Range("I2").Value = 19353.842
$ws.Range("J2").Value = 1910
$ws.Range("K2").Value = 19353.842
$ws.Range("L2").Value = 1910
$ws.Range("M2").Value = -19240.842
$ws.Range("N2").Value = -2136

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2828754
$ws.Range("I32").Value = 3400282.2
$ws.Range("J32").Value = 18740.834
$ws.Range("K32").Value = 3400282.2
$ws.Range("L32").Value = 18740.834
$ws.Range("M32").Value = -3399995.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2858.625
$ws.Range("I45").Value = 1417.7142
$ws.Range("J45").Value = 3979.3333
$ws.Range("K45").Value = 1417.7142
$ws.Range("L45").Value = 3979.3333
$ws.Range("M45").Value = -1040.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3243097.5
$ws.Range("I61").Value = 1737383
$ws.Range("J61").Value = 8405547
$ws.Range("K61").Value = 1737383
$ws.Range("L61").Value = 8405547
$ws.Range("M61").Value = -1737171
$ws.Range("N61").Value = -8405971

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 16975.137
$ws.Range("I116").Value = 19353.842
$ws.Range("J116").Value = 1910
$ws.Range("K116").Value = 19353.842
$ws.Range("L116").Value = 1910
$ws.Range("M116").Value = -17059.842
$ws.Range("N116").Value = -6498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 21988858
$ws.Range("I132").Value = 26201848
$ws.Range("J132").Value = 5558199.5
$ws.Range("K132").Value = 78605544
$ws.Range("L132").Value = 16674598.5
$ws.Range("M132").Value = -78603014
$ws.Range("N132").Value = -16679658.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3243097.5
$ws.Range("I136").Value = 1737383
$ws.Range("J136").Value = 8405547
$ws.Range("K136").Value = 5212149
$ws.Range("L136").Value = 25216641
$ws.Range("M136").Value = -5209599
$ws.Range("N136").Value = -25221741

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16975.137
$ws.Range("I3").Value = 19353.842
$ws.Range("J3").Value = 1910
$ws.Range("K3").Value = 19353.842
$ws.Range("L3").Value = 1910
$ws.Range("M3").Value = -19239.842
$ws.Range("N3").Value = -2138

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20011692
$ws.Range("I20").Value = 31257646
$ws.Range("J20").Value = 18885.334
$ws.Range("K20").Value = 31257646
$ws.Range("L20").Value = 18885.334
$ws.Range("M20").Value = -31257399
$ws.Range("N20").Value = -19379.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 100000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 100000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 100000
$ws.Range("N106").Value = -102524

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 15485993
$ws.Range("I134").Value = 17232580
$ws.Range("J134").Value = 6753059.5
$ws.Range("K134").Value = 51697740
$ws.Range("L134").Value = 20259178.5
$ws.Range("M134").Value = -51695205
$ws.Range("N134").Value = -20264248.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2693046
$ws.Range("I31").Value = 3969041.5
$ws.Range("J31").Value = 13455.5
$ws.Range("K31").Value = 3969041.5
$ws.Range("L31").Value = 13455.5
$ws.Range("M31").Value = -3968746.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2693046
$ws.Range("I34").Value = 3969041.5
$ws.Range("J34").Value = 13455.5
$ws.Range("K34").Value = 3969041.5
$ws.Range("L34").Value = 13455.5
$ws.Range("M34").Value = -3968839.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 856.25
$ws.Range("I35").Value = 641.6667
$ws.Range("J35").Value = 1500
$ws.Range("K35").Value = 641.6667
$ws.Range("L35").Value = 1500
$ws.Range("M35").Value = -347.6667
$ws.Range("N35").Value = -2088

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H119").Value = 27500
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 27500
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 27500
$ws.Range("N119").Value = -37176

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1409.1
$ws.Range("I132").Value = 1034.7916
$ws.Range("J132").Value = 2906.3333
$ws.Range("K132").Value = 3104.3748
$ws.Range("L132").Value = 8718.999899999999
$ws.Range("M132").Value = -574.3748000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2224933
$ws.Range("I134").Value = 2853
$ws.Range("J134").Value = 13335333
$ws.Range("K134").Value = 8559
$ws.Range("L134").Value = 40005999
$ws.Range("M134").Value = -6024
$ws.Range("N134").Value = -40011069

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7338415.5
$ws.Range("I70").Value = 4468764
$ws.Range("J70").Value = 10208067
$ws.Range("K70").Value = 4468764
$ws.Range("L70").Value = 10208067
$ws.Range("M70").Value = -4468494
$ws.Range("N70").Value = -10208607

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7338415.5
$ws.Range("I73").Value = 4468764
$ws.Range("J73").Value = 10208067
$ws.Range("K73").Value = 4468764
$ws.Range("L73").Value = 10208067
$ws.Range("M73").Value = -4467828
$ws.Range("N73").Value = -10209939

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H120").Value = 34450
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 34450
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 34450
$ws.Range("N120").Value = -44126

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 24150
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 24150
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 72450
$ws.Range("N134").Value = -77520

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 10272.053
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 10272.053
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 30816.159
$ws.Range("N136").Value = -35916.159

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 29605
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 29605
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 29605
$ws.Range("N121").Value = -33099
